$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.394.69"
$ws.Range("E2").Value = "  +0.83%  "

$ws.Range("D3").Value = "2.956.82"
$ws.Range("E3").Value = "  +2.83%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.66"
$ws.Range("E5").Value = "  +1.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "595.77"
$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +6.39%  "

$ws.Range("D10").Value = "2.955.62"
$ws.Range("E10").Value = "  +2.60%  "

$ws.Range("E11").Value = "  +12.23%  "

$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("D13").Value = "3.501.96"
$ws.Range("E13").Value = "  +2.63%  "

$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("E15").Value = "  +5.01%  "

$ws.Range("D16").Value = "76.374.34"
$ws.Range("E16").Value = "  +0.84%  "

$ws.Range("E17").Value = "  +0.83%  "

$ws.Range("D18").Value = "2.954.55"
$ws.Range("E18").Value = "  +1.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.68"
$ws.Range("E19").Value = "  +9.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.75"
$ws.Range("E20").Value = "  -1.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.85"
$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("E23").Value = "  +4.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.29"
$ws.Range("E24").Value = "  +1.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.31"
$ws.Range("E27").Value = "  +2.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.75"
$ws.Range("E28").Value = "  +0.42%  "

$ws.Range("E29").Value = "  +1.12%  "

$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.54"
$ws.Range("E31").Value = "  +10.13%  "

$ws.Range("E32").Value = "  -1.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "496.75"
$ws.Range("E33").Value = "  -1.94%  "

$ws.Range("E34").Value = "  +0.43%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.58"
$ws.Range("E36").Value = "  +0.67%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.40"
$ws.Range("E37").Value = "  +1.17%  "

$ws.Range("E39").Value = "  +20.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.96"
$ws.Range("E40").Value = "  +1.48%  "

$ws.Range("E41").Value = "  -1.86%  "

$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "180.84"
$ws.Range("E43").Value = "  -1.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.94"
$ws.Range("E44").Value = "  -1.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.65"
$ws.Range("E45").Value = "  -1.44%  "

$ws.Range("E46").Value = "  -0.95%  "

$ws.Range("E47").Value = "  -1.93%  "

$ws.Range("E48").Value = "  +2.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.90"
$ws.Range("E49").Value = "  +3.63%  "

$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.667"
$ws.Range("E51").Value = "  +0.04%  "

